# Implemented complexity 2 Instructions
# -Implemented: ADDIU, ANDI, LUI, ORI, SLL, SLT, SRA, SRAV, SRL, XORI.
# - ADD, ADDI, LW, SUB, SW not yet implemented (need memory and overflow function)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the Implementation? column (E) as implemented (1) for the newly
# implemented complexity-2 instructions.
$ws.Range("E4").Value = 1   # ADDIU
$ws.Range("E7").Value = 1   # ANDI
$ws.Range("E26").Value = 1  # LUI
$ws.Range("E37").Value = 1  # ORI
$ws.Range("E40").Value = 1  # SLL
$ws.Range("E42").Value = 1  # SLT

# LW still needs memory support - add a note explaining why.
$ws.Range("G27").Value = "memory not yet implemented"

# Update the view state to match where the user left off working.
$ws.Range("E43").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 12
$aw.ScrollColumn = 1

$excel.ActiveWindow.Left = 14880
